$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "41.528.37"
$ws.Range("E2").Value = "  +0.05%  "
Set-TextValue "D3" "2.470.81"
$ws.Range("E3").Value = "  -0.66%  "
Set-TextValue "D4" "0.998"
$ws.Range("E4").Value = "  -0.41%  "
Set-TextValue "D5" "314.68"
$ws.Range("E5").Value = "  -0.19%  "
Set-TextValue "D6" "92.03"
$ws.Range("E6").Value = "  -2.20%  "
$ws.Range("E7").Value = "  +0.08%  "
Set-TextValue "D8" "1.00"
$ws.Range("E8").Value = "  -0.36%  "
$ws.Range("E9").Value = "  +3.56%  "
Set-TextValue "D10" "32.54"
$ws.Range("E10").Value = "  -3.01%  "
Set-TextValue "D11" "0.0793"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +0.34%  "
Set-TextValue "D13" "2.851.24"
$ws.Range("E14").Value = "  -0.81%  "
Set-TextValue "D15" "15.98"
$ws.Range("E15").Value = "  +3.29%  "
Set-TextValue "D16" "2.471.03"
$ws.Range("E16").Value = "  -0.56%  "
Set-TextValue "D17" "0.778"
$ws.Range("E17").Value = "  -1.68%  "
Set-TextValue "D18" "41.538.36"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("E19").Value = "  +2.43%  "
Set-TextValue "D20" "0.0₃0943"
$ws.Range("E20").Value = "  +1.92%  "
Set-TextValue "D21" "71.02"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("E22").Value = "  -0.99%  "
Set-TextValue "D23" "238.09"
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +0.05%  "
Set-TextValue "D27" "24.80"
$ws.Range("E27").Value = "  +2.59%  "
$ws.Range("E28").Value = "  -1.42%  "
$ws.Range("E29").Value = "  -1.59%  "
Set-TextValue "D30" "35.52"
$ws.Range("E30").Value = "  -3.85%  "
Set-TextValue "D31" "156.02"
$ws.Range("E31").Value = "  +1.02%  "
$ws.Range("E32").Value = "  -0.94%  "
$ws.Range("E33").Value = "  +0.13%  "
Set-TextValue "D34" "0.0761"
$ws.Range("E34").Value = "  +0.84%  "
Set-TextValue "D35" "17.31"
$ws.Range("E35").Value = "  -3.47%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("E37").Value = "  -5.45%  "
$ws.Range("E38").Value = "  +3.28%  "
$ws.Range("E39").Value = "  +0.12%  "
Set-TextValue "D40" "1.80"
$ws.Range("E40").Value = "  -3.99%  "
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("E42").Value = "  -0.63%  "
Set-TextValue "D43" "1.946.97"
$ws.Range("E43").Value = "  -2.17%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D44" "0.0284"
$ws.Range("E44").Value = "  -0.84%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "18.92"
$ws.Range("E45").Value = "  -3.99%  "
$ws.Range("E46").Value = "  -3.58%  "
$ws.Range("E47").Value = "  +3.11%  "
Set-TextValue "D48" "2.709.21"
Set-TextValue "D49" "97.38"
$ws.Range("E49").Value = "  +0.27%  "
Set-TextValue "D50" "67.29"
$ws.Range("E50").Value = "  -3.10%  "
Set-TextValue "D51" "52.39"
$ws.Range("E51").Value = "  +3.11%  "
